# Cambio transacción ZMM023, si no encuentra tickets de registro, avanzar a
# procesar tickets en estado REV.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the execution date value (A2) to the new date.
$ws.Range("A2").Value = "18/01/2023"

# Clear out the stale historical execution dates (A3:A8) - only the
# current date should remain in the sheet now.
$ws.Range("A3:A8").ClearContents()

# Move the active selection from B7 to D7.
$ws.Range("D7").Select()

$wb.Save()
